$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Remove the two bullet paragraphs that are being dropped:
#   "Avanço da tecnologia provém novas ferramentas"
#   "Ferramentas aplicáveis em distintos modelos de negócios"
# Delete paragraph 2 first so paragraph 1's index stays valid.
$tr.Paragraphs(2, 1).Delete()
$tr.Paragraphs(1, 1).Delete()

# The remaining first paragraph is now "Aplicação nos contextos:" -
# split it into two runs: "Aplicação " and "nos contextos:"
$para = $tr.Paragraphs(1, 1)
$splitPoint = "Aplicação ".Length
$firstRun = $para.Characters(1, $splitPoint)
$firstRun.Text = "Aplicação "

# Reset autofit so the (now shorter) text isn't shrunk anymore.
$tf.AutoSize = 2
